$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per diff. D-column (Price) values are numeric-looking text
# (e.g. "1.000", "30.168.32") that must stay plain text, exactly as stored
# in the original workbook (inline/shared string, not a Number cell).
# Force text interpretation: set NumberFormat to "@" before assigning the
# value, then reset the cell style to "Normal" so no stray custom style/
# number-format lingers on the cell afterwards.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.168.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.905.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  -3.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4026"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08275"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.109"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.401"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.889.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.307"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.92%  "
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06460"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.945"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.200.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.200"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.59%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.117.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.337"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.125"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1043"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.983"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.737"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02456"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.360"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06454"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2156"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.52%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.670"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.184"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6388"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.218"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.188"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.74%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5986"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.645"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.216"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.14%  "
